# RegistrationData.xlsx edit
# - Row 7 (C7): email value changed from "demo1@example.com" to "demo4@example.com"
# - Row 9 (D9): phone number value changed from numeric 12 to text "1"
# - Row 12 (D12): phone number value changed from numeric 2.12343456521312E+35
#                 to text "22222222222222222222222222222222222"
# - Active selection moved from D16 to D14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hyperlinked email text in row 7, column C
$ws.Range("C7").Value = "demo4@example.com"

# D12 must be entered first so its new shared-string slot precedes the one
# created for D9 (matches the original authoring order).
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22222222222222222222222222222222222"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1"

# Move the active cell selection to D14
$ws.Range("D14").Select()
